# create biaya pendaftaran, spp, bpp
#
# Adds three new rows of data (Pendaftaran / BPP / SPP) into column B of
# rows 5-7, matching the formatting already used by the rows above them:
#   - column B keeps the "label" style (border, no fill)
#   - columns C:F get the green-fill/border style used by C2:F4
#   - column G gets the black-fill/border style used by G2
# Also normalizes the G1 header cell's fill (it was redundantly using a
# "fill applied but no actual fill colour" style) to the same plain
# centered/bordered style already used by the rest of the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize the G1 header style so it matches A1:F1 exactly ---------
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats

# --- New data rows: labels ---------------------------------------------
$ws.Range("B5").Value = "Pendaftaran"
$ws.Range("B6").Value = "BPP"
$ws.Range("B7").Value = "SPP"

# --- Copy formatting for C:F (green fill) from the existing C2:F2 ------
$ws.Range("C2:F2").Copy()
$ws.Range("C5:F5").PasteSpecial(-4122)
$ws.Range("C6:F6").PasteSpecial(-4122)
$ws.Range("C7:F7").PasteSpecial(-4122)

# --- Copy formatting for G (black fill) from the existing G2 -----------
$ws.Range("G2").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("G7").PasteSpecial(-4122)

# Clear clipboard marquee / restore cursor position like the author left it
$excel.CutCopyMode = $false
[void]$ws.Range("H14").Select()

Write-Output "Added Pendaftaran/BPP/SPP rows"
